$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "246.95"
Set-TextValue $ws.Range("G2") "16"
Set-TextValue $ws.Range("D3") "29.37"
Set-TextValue $ws.Range("E3") "7.68%"
Set-TextValue $ws.Range("G3") "16"
Set-TextValue $ws.Range("D4") "5.189"
Set-TextValue $ws.Range("E4") "2.70%"
Set-TextValue $ws.Range("G4") "16"
Set-TextValue $ws.Range("D5") "0.05709"
Set-TextValue $ws.Range("E5") "0.35%"
Set-TextValue $ws.Range("G5") "16"
Set-TextValue $ws.Range("D6") "6.581"
Set-TextValue $ws.Range("E6") "1.20%"
Set-TextValue $ws.Range("G6") "16"
Set-TextValue $ws.Range("B7") "MXToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D7") "0.8582"
Set-TextValue $ws.Range("E7") "4.61%"
Set-TextValue $ws.Range("G7") "16"
Set-TextValue $ws.Range("B8") "FTXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D8") "0.8770"
Set-TextValue $ws.Range("E8") "4.47%"
Set-TextValue $ws.Range("G8") "16"
Set-TextValue $ws.Range("B9") "WazirX"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D9") "0.1369"
Set-TextValue $ws.Range("E9") "3.05%"
Set-TextValue $ws.Range("G9") "16"
Set-TextValue $ws.Range("B10") "MandalaExchangeToken"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D10") "0.07067"
Set-TextValue $ws.Range("E10") "2.00%"
Set-TextValue $ws.Range("G10") "16"
Set-TextValue $ws.Range("B11") "BitrueCoin"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D11") "0.02870"
Set-TextValue $ws.Range("E11") "1.16%"
Set-TextValue $ws.Range("G11") "16"
Set-TextValue $ws.Range("B12") "BitMartToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D12") "0.09386"
Set-TextValue $ws.Range("E12") "-0.16%"
Set-TextValue $ws.Range("G12") "16"
Set-TextValue $ws.Range("B13") "BitForexToken"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D13") "0.001526"
Set-TextValue $ws.Range("E13") "1.24%"
Set-TextValue $ws.Range("G13") "16"
Set-TextValue $ws.Range("B14") "CoinExToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D14") "0.04187"
Set-TextValue $ws.Range("E14") "2.12%"
Set-TextValue $ws.Range("G14") "16"
Set-TextValue $ws.Range("B15") "One"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D15") "0.0006007"
Set-TextValue $ws.Range("E15") "-94.01%"
Set-TextValue $ws.Range("G15") "16"
Set-TextValue $ws.Range("D16") "0.006142"
Set-TextValue $ws.Range("E16") "-0.85%"
Set-TextValue $ws.Range("G16") "16"
Set-TextValue $ws.Range("E17") "3,769.23%"
Set-TextValue $ws.Range("G17") "16"
Set-TextValue $ws.Range("D18") "3.481"
Set-TextValue $ws.Range("E18") "-0.82%"
Set-TextValue $ws.Range("G18") "16"
Set-TextValue $ws.Range("B19") "GateToken"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D19") "3.077"
Set-TextValue $ws.Range("E19") "2.48%"
Set-TextValue $ws.Range("G19") "16"
Set-TextValue $ws.Range("B20") "BTSEToken"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D20") "2.278"
Set-TextValue $ws.Range("E20") "-1.32%"
Set-TextValue $ws.Range("G20") "16"
Set-TextValue $ws.Range("D21") "0.3174"
Set-TextValue $ws.Range("E21") "0.68%"
Set-TextValue $ws.Range("G21") "16"
Set-TextValue $ws.Range("D22") "0.03311"
Set-TextValue $ws.Range("E22") "4.43%"
Set-TextValue $ws.Range("G22") "16"
Set-TextValue $ws.Range("D23") "0.1301"
Set-TextValue $ws.Range("E23") "3.60%"
Set-TextValue $ws.Range("G23") "16"
Set-TextValue $ws.Range("D24") "3.463"
Set-TextValue $ws.Range("E24") "-2.97%"
Set-TextValue $ws.Range("G24") "16"
Set-TextValue $ws.Range("D25") "0.1379"
Set-TextValue $ws.Range("E25") "0.47%"
Set-TextValue $ws.Range("G25") "16"
Set-TextValue $ws.Range("D26") "0.005048"
Set-TextValue $ws.Range("E26") "27.88%"
Set-TextValue $ws.Range("G26") "16"
Set-TextValue $ws.Range("D27") "0.001221"
Set-TextValue $ws.Range("E27") "0.45%"
Set-TextValue $ws.Range("G27") "16"
Set-TextValue $ws.Range("E28") "23.51%"
Set-TextValue $ws.Range("G28") "16"
Set-TextValue $ws.Range("G29") "16"
Set-TextValue $ws.Range("G30") "16"
Set-TextValue $ws.Range("G31") "16"
Set-TextValue $ws.Range("G32") "16"
Set-TextValue $ws.Range("G33") "16"
Set-TextValue $ws.Range("G34") "16"
Set-TextValue $ws.Range("G35") "16"
Set-TextValue $ws.Range("G36") "16"
Set-TextValue $ws.Range("G37") "16"
Set-TextValue $ws.Range("G38") "16"
Set-TextValue $ws.Range("G39") "16"
Set-TextValue $ws.Range("D40") "0.03752"
Set-TextValue $ws.Range("E40") "1.45%"
Set-TextValue $ws.Range("G40") "16"
Set-TextValue $ws.Range("D41") "0.005697"
Set-TextValue $ws.Range("E41") "-0.67%"
Set-TextValue $ws.Range("G41") "16"
Set-TextValue $ws.Range("D42") "0.1072"
Set-TextValue $ws.Range("E42") "1.96%"
Set-TextValue $ws.Range("G42") "16"
Set-TextValue $ws.Range("D43") "0.001999"
Set-TextValue $ws.Range("E43") "11.16%"
Set-TextValue $ws.Range("G43") "16"
Set-TextValue $ws.Range("D44") "0.01025"
Set-TextValue $ws.Range("E44") "6.83%"
Set-TextValue $ws.Range("G44") "16"
Set-TextValue $ws.Range("D45") "0.00005158"
Set-TextValue $ws.Range("E45") "-0.96%"
Set-TextValue $ws.Range("G45") "16"
Set-TextValue $ws.Range("D46") "0.00000000750"
Set-TextValue $ws.Range("E46") "0.04%"
Set-TextValue $ws.Range("G46") "16"
Set-TextValue $ws.Range("D47") "0.06996"
Set-TextValue $ws.Range("E47") "-31.01%"
Set-TextValue $ws.Range("G47") "16"
Set-TextValue $ws.Range("D48") "0.002578"
Set-TextValue $ws.Range("E48") "0.36%"
Set-TextValue $ws.Range("G48") "16"
Set-TextValue $ws.Range("D49") "0.00002099"
Set-TextValue $ws.Range("E49") "0.04%"
Set-TextValue $ws.Range("G49") "16"
Set-TextValue $ws.Range("D50") "0.0001999"
Set-TextValue $ws.Range("E50") "0.04%"
Set-TextValue $ws.Range("G50") "16"
Set-TextValue $ws.Range("G51") "16"
